$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": update changed cells (rows 2-44), then add new rows 45-46 ---
# Row 2
$ws1.Cells.Item(2,4).Value2 = 3229.89
$ws1.Cells.Item(2,5).Value2 = 99.66
# Row 3
$ws1.Cells.Item(3,4).Value2 = 2549.29
$ws1.Cells.Item(3,5).Value2 = 637.7
# Row 5
$ws1.Cells.Item(5,4).Value2 = 2285
# Row 6
$ws1.Cells.Item(6,4).Value2 = 2185
$ws1.Cells.Item(6,5).Value2 = 545
# Row 7
$ws1.Cells.Item(7,4).Value2 = 2185
$ws1.Cells.Item(7,5).Value2 = 580
# Row 8
$ws1.Cells.Item(8,1).Value2 = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(8,3).Value2 = 4
$ws1.Cells.Item(8,4).Value2 = 1462.23
$ws1.Cells.Item(8,5).Value2 = 371.75
# Row 9
$ws1.Cells.Item(9,1).Value2 = "BRVM - TRANSPORT"
$ws1.Cells.Item(9,4).Value2 = 1396.44
$ws1.Cells.Item(9,5).Value2 = 353.7
# Row 10
$ws1.Cells.Item(10,1).Value2 = "BRVM - AGRICULTURE"
$ws1.Cells.Item(10,4).Value2 = 1290.26
$ws1.Cells.Item(10,5).Value2 = 326.77
# Row 11
$ws1.Cells.Item(11,1).Value2 = "CFAO MOTORS CI"
$ws1.Cells.Item(11,3).Value2 = 2
$ws1.Cells.Item(11,4).Value2 = 1265
$ws1.Cells.Item(11,5).Value2 = 635
# Row 12
$ws1.Cells.Item(12,4).Value2 = 1053.58
$ws1.Cells.Item(12,5).Value2 = 262.08
# Row 13
$ws1.Cells.Item(13,4).Value2 = 869.48
$ws1.Cells.Item(13,5).Value2 = 217.18
# Row 14
$ws1.Cells.Item(14,4).Value2 = 757.59
$ws1.Cells.Item(14,5).Value2 = 189.81
# Row 15
$ws1.Cells.Item(15,4).Value2 = 547.04
$ws1.Cells.Item(15,5).Value2 = 137.18
# Row 16
$ws1.Cells.Item(16,4).Value2 = 516.36
$ws1.Cells.Item(16,5).Value2 = 129.75
# Row 17
$ws1.Cells.Item(17,4).Value2 = 487.44
$ws1.Cells.Item(17,5).Value2 = 122.1
# Row 18
$ws1.Cells.Item(18,4).Value2 = 479.04
$ws1.Cells.Item(18,5).Value2 = 119.99
# Row 19
$ws1.Cells.Item(19,4).Value2 = 436.67
$ws1.Cells.Item(19,5).Value2 = 110.63
# Row 20
$ws1.Cells.Item(20,4).Value2 = 420.68
$ws1.Cells.Item(20,5).Value2 = 106.5
# Row 21
$ws1.Cells.Item(21,4).Value2 = 369.41
$ws1.Cells.Item(21,5).Value2 = 93.31999999999999
# Row 22
$ws1.Cells.Item(22,2).Value2 = 2
$ws1.Cells.Item(22,4).Value2 = 10.47
$ws1.Cells.Item(22,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(22,7).Value2 = "➖ Neutre"
# Row 23
$ws1.Cells.Item(23,2).Value2 = 1
$ws1.Cells.Item(23,4).Value2 = 5.73
# Row 24
$ws1.Cells.Item(24,1).Value2 = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(24,2).Value2 = 1
$ws1.Cells.Item(24,3).Value2 = 0
$ws1.Cells.Item(24,4).Value2 = 5.45
$ws1.Cells.Item(24,5).Value2 = 5.45
$ws1.Cells.Item(24,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(24,7).Value2 = "➖ Neutre"
# Row 25
$ws1.Cells.Item(25,3).Value2 = 1
$ws1.Cells.Item(25,4).Value2 = 3.96
$ws1.Cells.Item(25,5).Value2 = -3.45
$ws1.Cells.Item(25,7).Value2 = "👀 À surveiller"
# Row 26
$ws1.Cells.Item(26,1).Value2 = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(26,2).Value2 = 1
$ws1.Cells.Item(26,3).Value2 = 0
$ws1.Cells.Item(26,4).Value2 = 3.54
$ws1.Cells.Item(26,5).Value2 = 3.54
$ws1.Cells.Item(26,7).Value2 = "➖ Neutre"
# Row 27
$ws1.Cells.Item(27,1).Value2 = "PALM CI (PALC)"
$ws1.Cells.Item(27,2).Value2 = 1
$ws1.Cells.Item(27,3).Value2 = 0
$ws1.Cells.Item(27,4).Value2 = 3.41
$ws1.Cells.Item(27,5).Value2 = 3.41
$ws1.Cells.Item(27,7).Value2 = "➖ Neutre"
# Row 28
$ws1.Cells.Item(28,1).Value2 = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(28,3).Value2 = 1
$ws1.Cells.Item(28,4).Value2 = 3.1
$ws1.Cells.Item(28,5).Value2 = -0.6
$ws1.Cells.Item(28,7).Value2 = "👀 À surveiller"
# Row 30
$ws1.Cells.Item(30,1).Value2 = "CFAO MOTORS CI (CFAC)"
$ws1.Cells.Item(30,3).Value2 = 0
$ws1.Cells.Item(30,4).Value2 = 3.05
$ws1.Cells.Item(30,5).Value2 = 3.05
$ws1.Cells.Item(30,7).Value2 = "➖ Neutre"
# Row 31
$ws1.Cells.Item(31,1).Value2 = "ORANGE COTE D'IVOIRE (ORAC)"
$ws1.Cells.Item(31,3).Value2 = 0
$ws1.Cells.Item(31,4).Value2 = 2.93
$ws1.Cells.Item(31,5).Value2 = 2.93
$ws1.Cells.Item(31,7).Value2 = "➖ Neutre"
# Row 32
$ws1.Cells.Item(32,1).Value2 = "SICABLE CI (CABC)"
$ws1.Cells.Item(32,3).Value2 = 0
$ws1.Cells.Item(32,4).Value2 = 2.69
$ws1.Cells.Item(32,5).Value2 = 2.69
$ws1.Cells.Item(32,7).Value2 = "➖ Neutre"
# Row 33
$ws1.Cells.Item(33,1).Value2 = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(33,2).Value2 = 1
$ws1.Cells.Item(33,3).Value2 = 1
$ws1.Cells.Item(33,4).Value2 = 2.02
$ws1.Cells.Item(33,5).Value2 = -2.02
$ws1.Cells.Item(33,7).Value2 = "👀 À surveiller"
# Row 34
$ws1.Cells.Item(34,1).Value2 = "ONATEL BF (ONTBF)"
$ws1.Cells.Item(34,2).Value2 = 1
$ws1.Cells.Item(34,3).Value2 = 0
$ws1.Cells.Item(34,4).Value2 = 1.96
$ws1.Cells.Item(34,5).Value2 = 1.96
# Row 35
$ws1.Cells.Item(35,1).Value2 = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(35,2).Value2 = 2
$ws1.Cells.Item(35,4).Value2 = 0.79
$ws1.Cells.Item(35,5).Value2 = -5.88
# Row 36
$ws1.Cells.Item(36,1).Value2 = "TOTAL"
$ws1.Cells.Item(36,3).Value2 = 4
$ws1.Cells.Item(36,4).Value2 = 0
$ws1.Cells.Item(36,5).Value2 = 0
# Row 37
$ws1.Cells.Item(37,1).Value2 = "BICI CI (BICC)"
$ws1.Cells.Item(37,4).Value2 = -0.6
$ws1.Cells.Item(37,5).Value2 = -0.6
# Row 38
$ws1.Cells.Item(38,1).Value2 = "SAPH CI (SPHC)"
$ws1.Cells.Item(38,2).Value2 = 1
$ws1.Cells.Item(38,3).Value2 = 2
$ws1.Cells.Item(38,4).Value2 = -1.42
$ws1.Cells.Item(38,5).Value2 = 4.03
$ws1.Cells.Item(38,7).Value2 = "👀 À surveiller"
# Row 39
$ws1.Cells.Item(39,1).Value2 = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Cells.Item(39,4).Value2 = -1.42
$ws1.Cells.Item(39,5).Value2 = -1.42
# Row 40
$ws1.Cells.Item(40,1).Value2 = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(40,4).Value2 = -1.75
$ws1.Cells.Item(40,5).Value2 = -1.75
# Row 41
$ws1.Cells.Item(41,1).Value2 = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(41,4).Value2 = -1.87
$ws1.Cells.Item(41,5).Value2 = -1.87
# Row 42
$ws1.Cells.Item(42,1).Value2 = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Cells.Item(42,4).Value2 = -2
$ws1.Cells.Item(42,5).Value2 = -2
# Row 43
$ws1.Cells.Item(43,1).Value2 = "SOGB CI (SOGC)"
$ws1.Cells.Item(43,2).Value2 = 0
$ws1.Cells.Item(43,4).Value2 = -3.23
$ws1.Cells.Item(43,5).Value2 = -3.23
$ws1.Cells.Item(43,7).Value2 = "➖ Neutre"
# Row 44
$ws1.Cells.Item(44,1).Value2 = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(44,2).Value2 = 1
$ws1.Cells.Item(44,3).Value2 = 1
$ws1.Cells.Item(44,4).Value2 = -3.46
$ws1.Cells.Item(44,5).Value2 = 3.09
$ws1.Cells.Item(44,7).Value2 = "👀 À surveiller"
# Row 45 (new)
$ws1.Cells.Item(45,1).Value2 = "SAFCA CI (SAFC)"
$ws1.Cells.Item(45,2).Value2 = 2
$ws1.Cells.Item(45,3).Value2 = 2
$ws1.Cells.Item(45,4).Value2 = -4.14
$ws1.Cells.Item(45,5).Value2 = -7.33
$ws1.Cells.Item(45,6).Value2 = "🟡 Observer"
$ws1.Cells.Item(45,7).Value2 = "👀 À surveiller"
# Row 46 (new)
$ws1.Cells.Item(46,1).Value2 = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(46,2).Value2 = 0
$ws1.Cells.Item(46,3).Value2 = 3
$ws1.Cells.Item(46,4).Value2 = -7.17
$ws1.Cells.Item(46,5).Value2 = -1.97
$ws1.Cells.Item(46,6).Value2 = "🔴 Vente"
$ws1.Cells.Item(46,7).Value2 = "⚠️ Risque de décrochage"

# --- Sheet "Top_YTD": update changed cells ---
# Row 2
$ws2.Cells.Item(2,2).Value2 = 6803240.42
# Row 3
$ws2.Cells.Item(3,2).Value2 = 295442.34
# Row 5
$ws2.Cells.Item(5,2).Value2 = 202778
# Row 6
$ws2.Cells.Item(6,2).Value2 = 174308
# Row 7
$ws2.Cells.Item(7,2).Value2 = 174001.76
# Row 8
$ws2.Cells.Item(8,2).Value2 = 46866.18
# Row 9
$ws2.Cells.Item(9,2).Value2 = 40577.85
# Row 10
$ws2.Cells.Item(10,1).Value2 = "BRVM - AGRICULTURE"
$ws2.Cells.Item(10,2).Value2 = 31778.83
# Row 11
$ws2.Cells.Item(11,1).Value2 = "BRVM - INDUSTRIE"
$ws2.Cells.Item(11,2).Value2 = 17338.6
